$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell $ws 'D2' '69.286.45'
Set-TextCell $ws 'E2' '  +4.83%  '
Set-TextCell $ws 'D3' '3.583.53'
Set-TextCell $ws 'E3' '  +16.94%  '
Set-TextCell $ws 'E4' '  -0.05%  '
Set-TextCell $ws 'D5' '590.54'
Set-TextCell $ws 'E5' '  +3.14%  '
Set-TextCell $ws 'D6' '184.81'
Set-TextCell $ws 'D7' '3.585.45'
Set-TextCell $ws 'E7' '  +17.09%  '
Set-TextCell $ws 'D8' '0.999'
Set-TextCell $ws 'E8' '  -0.14%  '
Set-TextCell $ws 'D9' '0.532'
Set-TextCell $ws 'E9' '  +4.67%  '
Set-TextCell $ws 'B10' 'Dogecoin'
Set-TextCell $ws 'C10' 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell $ws 'D10' '0.159'
Set-TextCell $ws 'E10' '  +7.24%  '
Set-TextCell $ws 'B11' 'Toncoin'
Set-TextCell $ws 'C11' 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell $ws 'D11' '6.59'
Set-TextCell $ws 'E11' '  +4.74%  '
Set-TextCell $ws 'D12' '0.493'
Set-TextCell $ws 'E12' '  +5.78%  '
Set-TextCell $ws 'B13' 'Avalanche'
Set-TextCell $ws 'C13' 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
Set-TextCell $ws 'D13' '38.60'
Set-TextCell $ws 'E13' '  +8.40%  '
Set-TextCell $ws 'B14' 'ShibaInu'
Set-TextCell $ws 'C14' 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
Set-TextCell $ws 'D14' '0.0000251'
Set-TextCell $ws 'E14' '  +5.51%  '
Set-TextCell $ws 'D15' '4.162.25'
Set-TextCell $ws 'E15' '  +16.43%  '
Set-TextCell $ws 'D16' '69.496.13'
Set-TextCell $ws 'E16' '  +5.29%  '
Set-TextCell $ws 'D17' '3.565.86'
Set-TextCell $ws 'E17' '  +16.46%  '
Set-TextCell $ws 'E18' '  +1.56%  '
Set-TextCell $ws 'D19' '7.49'
Set-TextCell $ws 'E19' '  +8.21%  '
Set-TextCell $ws 'D20' '16.95'
Set-TextCell $ws 'E20' '  +3.25%  '
Set-TextCell $ws 'D21' '509.49'
Set-TextCell $ws 'E21' '  +5.51%  '
Set-TextCell $ws 'D22' '9.24'
Set-TextCell $ws 'E22' '  +21.26%  '
Set-TextCell $ws 'D23' '0.741'
Set-TextCell $ws 'E23' '  +8.44%  '
Set-TextCell $ws 'D24' '86.71'
Set-TextCell $ws 'E24' '  +5.45%  '
Set-TextCell $ws 'D25' '13.41'
Set-TextCell $ws 'E25' '  +6.68%  '
Set-TextCell $ws 'D26' '2.38'
Set-TextCell $ws 'E26' '  +8.41%  '
Set-TextCell $ws 'D27' '10.75'
Set-TextCell $ws 'E27' '  +6.11%  '
Set-TextCell $ws 'E28' '  +0.07%  '
Set-TextCell $ws 'D29' '2.53'
Set-TextCell $ws 'E29' '  +12.90%  '
Set-TextCell $ws 'D30' '8.07'
Set-TextCell $ws 'E30' '  +2.76%  '
Set-TextCell $ws 'D31' '31.84'
Set-TextCell $ws 'E31' '  +15.57%  '
Set-TextCell $ws 'E32' '  +21.66%  '
Set-TextCell $ws 'D33' '2.73'
Set-TextCell $ws 'E33' '  +5.55%  '
Set-TextCell $ws 'D34' '0.117'
Set-TextCell $ws 'E34' '  +5.98%  '
Set-TextCell $ws 'D35' '0.997'
Set-TextCell $ws 'E35' '  -0.13%  '
Set-TextCell $ws 'D36' '6.12'
Set-TextCell $ws 'E36' '  +10.76%  '
Set-TextCell $ws 'D37' '1.02'
Set-TextCell $ws 'E37' '  +8.90%  '
Set-TextCell $ws 'D38' '0.332'
Set-TextCell $ws 'E38' '  +11.38%  '
Set-TextCell $ws 'B39' 'Arweave'
Set-TextCell $ws 'C39' 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
Set-TextCell $ws 'D39' '46.84'
Set-TextCell $ws 'E39' '  -0.83%  '
Set-TextCell $ws 'B40' 'Stacks'
Set-TextCell $ws 'C40' 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextCell $ws 'D40' '2.09'
Set-TextCell $ws 'E40' '  +7.72%  '
Set-TextCell $ws 'D41' '50.64'
Set-TextCell $ws 'E41' '  +3.40%  '
Set-TextCell $ws 'E42' '  +4.75%  '
Set-TextCell $ws 'D43' '8.78'
Set-TextCell $ws 'E43' '  +7.05%  '
Set-TextCell $ws 'D44' '3.069.52'
Set-TextCell $ws 'E44' '  +11.22%  '
Set-TextCell $ws 'D45' '2.82'
Set-TextCell $ws 'E45' '  +12.42%  '
Set-TextCell $ws 'D46' '401.62'
Set-TextCell $ws 'E46' '  +11.22%  '
Set-TextCell $ws 'D47' '0.0364'
Set-TextCell $ws 'E47' '  +6.63%  '
Set-TextCell $ws 'D48' '27.61'
Set-TextCell $ws 'E48' '  +14.43%  '
Set-TextCell $ws 'D49' '135.33'
Set-TextCell $ws 'E49' '  +0.69%  '
Set-TextCell $ws 'E50' '  +0.06%  '
Set-TextCell $ws 'D51' '2.45'
Set-TextCell $ws 'E51' '  +14.95%  '
